# Add a new "2022-Q1" sheet (fund holdings detail) between "2021-Q3" and "总计",
# and add a new summary row to "总计" for the 2022-Q1 quarter.

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# 1) Build the new "2022-Q1" sheet by duplicating "总计" (so it inherits the
#    same header/column-A cell styling) and placing the copy right before it.
#    NOTE: after .Copy(Before:=totalSheet), the $totalSheet handle itself
#    tracks the freshly-created copy (it now sits at the original position),
#    so both sheets must be re-acquired by name afterwards.
# ---------------------------------------------------------------------------
$totalSheet.Copy($totalSheet, $null)

$newSheet = $wb.Worksheets.Item("总计 (2)")
$newSheet.Name = "2022-Q1"
$totalSheet = $wb.Worksheets.Item("总计")

# Extend the inherited styling (same-sheet copy keeps the style index) from
# the existing styled cells out to the extra columns/rows we need.
$newSheet.Range("D1").Copy()
$newSheet.Range("E1:H1").PasteSpecial(-4122)
$newSheet.Range("A2").Copy()
$newSheet.Range("A3:A6").PasteSpecial(-4122)

# Clear the leftover "总计"-style data row (2021-Q3 / 1 / 0.26) that came
# along with the duplicated sheet before writing the fund rows.
$newSheet.Range("B2:D2").ClearContents()

# Header row
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Make sure the numeric-looking text columns (B, D:G) stay TEXT, not numbers
# (fund codes have leading zeros; percentages/amounts are stored as text in
# the source data). NumberFormat="@" forces text-on-assignment; ClearFormats
# afterwards drops the residual style index again (the source file has no
# explicit style on these cells) while leaving the text type/value intact.
$newSheet.Range("B2:G6").NumberFormat = "@"

# Row index column (A) + fund rows (B:H)
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "501079"
$newSheet.Cells.Item(2,3).Value = "大成科创主题 3 年封闭运作灵活配置混合"
$newSheet.Cells.Item(2,4).Value = "17.69"
$newSheet.Cells.Item(2,5).Value = "79.13"
$newSheet.Cells.Item(2,6).Value = "4.78"
$newSheet.Cells.Item(2,7).Value = "0.8456"
$newSheet.Cells.Item(2,8).Value = 2

$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "012473"
$newSheet.Cells.Item(3,3).Value = "大成成长回报六个月持有期混合型证券投资基金A"
$newSheet.Cells.Item(3,4).Value = "8.97"
$newSheet.Cells.Item(3,5).Value = "71.30"
$newSheet.Cells.Item(3,6).Value = "4.82"
$newSheet.Cells.Item(3,7).Value = "0.4324"
$newSheet.Cells.Item(3,8).Value = 1

$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).Value = "010371"
$newSheet.Cells.Item(4,3).Value = "大成成长进取混合A"
$newSheet.Cells.Item(4,4).Value = "5.55"
$newSheet.Cells.Item(4,5).Value = "80.17"
$newSheet.Cells.Item(4,6).Value = "4.84"
$newSheet.Cells.Item(4,7).Value = "0.2686"
$newSheet.Cells.Item(4,8).Value = 3

$newSheet.Cells.Item(5,1).Value = 3
$newSheet.Cells.Item(5,2).Value = "010372"
$newSheet.Cells.Item(5,3).Value = "大成成长进取混合C"
$newSheet.Cells.Item(5,4).Value = "1.71"
$newSheet.Cells.Item(5,5).Value = "80.17"
$newSheet.Cells.Item(5,6).Value = "4.84"
$newSheet.Cells.Item(5,7).Value = "0.0828"
$newSheet.Cells.Item(5,8).Value = 3

$newSheet.Cells.Item(6,1).Value = 4
$newSheet.Cells.Item(6,2).Value = "012474"
$newSheet.Cells.Item(6,3).Value = "大成成长回报六个月持有期混合型证券投资基金C"
$newSheet.Cells.Item(6,4).Value = "0.43"
$newSheet.Cells.Item(6,5).Value = "71.30"
$newSheet.Cells.Item(6,6).Value = "4.82"
$newSheet.Cells.Item(6,7).Value = "0.0207"
$newSheet.Cells.Item(6,8).Value = 1

# Drop the residual text-number-format style picked up above; the source
# data has no explicit style on these cells (only plain inline strings).
$newSheet.Range("B2:G6").ClearFormats()

# ---------------------------------------------------------------------------
# 2) Insert the 2022-Q1 summary row at the top of "总计" (pushing the
#    existing 2021-Q3 row down), keeping the A-column index/style.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("B1").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 5
$totalSheet.Cells.Item(2,4).Value = 1.65

$totalSheet.Cells.Item(3,1).Value = 1

# Restore the original active sheet/selection ("2021-Q3" was active before
# this edit, and the diff leaves bookViews/activeTab untouched).
$sheet1 = $wb.Worksheets.Item("2021-Q3")
$sheet1.Activate() | Out-Null
$sheet1.Range("A1").Select() | Out-Null

